$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.233789666666667
$ws.Range("H2").Value = 3.701369
$ws.Range("I2").Value = 0.0001664233864291757
$ws.Range("J2").Value = 0.0001664233864291757
$ws.Range("Q2").Value = 0.00874839126311111
$ws.Range("R2").Value = 0.07873552136799999
$ws.Range("S2").Value = 0.0001664233864291757
$ws.Range("T2").Value = 0.0001664233864291757

# Row 3 updates
$ws.Range("I3").Value = 0.9827534361704352
$ws.Range("J3").Value = 0.9827534361704352
$ws.Range("S3").Value = 0.9827534361704352
$ws.Range("T3").Value = 0.9827534361704352

# Row 4 updates
$ws.Range("I4").Value = 0.01708014044313564
$ws.Range("J4").Value = 0.01708014044313564
$ws.Range("S4").Value = 0.01708014044313564
$ws.Range("T4").Value = 0.01708014044313564
